$wb = $excel.ActiveWorkbook

# Sheet ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 100001600
$ws.Range("I100").Value = 1995
$ws.Range("J100").Value = 500000000
$ws.Range("K100").Value = 1995
$ws.Range("L100").Value = 500000000
$ws.Range("M100").Value = -1454
$ws.Range("N100").Value = -500001082

# Sheet ALC row 108
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 34248.832
$ws.Range("J108").Value = 34248.832
$ws.Range("L108").Value = 34248.832
$ws.Range("N108").Value = -41928.832

# Sheet ALC row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 43494.285
$ws.Range("J130").Value = 43494.285
$ws.Range("L130").Value = 43494.285
$ws.Range("N130").Value = -53534.285

# Sheet ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 49722.832
$ws.Range("J133").Value = 49722.832
$ws.Range("L133").Value = 49722.832
$ws.Range("N133").Value = -59842.832

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7898.64
$ws.Range("I32").Value = 4900.2393
$ws.Range("J32").Value = 42380.25
$ws.Range("K32").Value = 4900.2393
$ws.Range("L32").Value = 42380.25
$ws.Range("M32").Value = -4613.2393
$ws.Range("N32").Value = -42954.25

# Sheet ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 32244
$ws.Range("J80").Value = 38030
$ws.Range("L80").Value = 38030
$ws.Range("N80").Value = -40026

# Sheet ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 32244
$ws.Range("J83").Value = 38030
$ws.Range("L83").Value = 114090
$ws.Range("N83").Value = -124074

# Sheet ARM row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 35200.668
$ws.Range("J101").Value = 35200.668
$ws.Range("L101").Value = 35200.668
$ws.Range("N101").Value = -41690.668

# Sheet ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 26844.25
$ws.Range("J109").Value = 26844.25
$ws.Range("L109").Value = 26844.25
$ws.Range("N109").Value = -29618.25

# Sheet ARM row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 10641.333
$ws.Range("J112").Value = 10641.333
$ws.Range("L112").Value = 10641.333
$ws.Range("N112").Value = -13595.333

# Sheet ARM row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 225035170
$ws.Range("J125").Value = 225035170
$ws.Range("L125").Value = 225035170
$ws.Range("N125").Value = -225045010

# Sheet ARM row 128
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 42500
$ws.Range("J128").Value = 42500
$ws.Range("L128").Value = 42500
$ws.Range("N128").Value = -52460

# Sheet ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 36701.75
$ws.Range("J135").Value = 36701.75
$ws.Range("L135").Value = 36701.75
$ws.Range("N135").Value = -46841.75

# Sheet BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 34246.332
$ws.Range("J82").Value = 35297.6
$ws.Range("L82").Value = 35297.6
$ws.Range("N82").Value = -36063.6

# Sheet BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 34246.332
$ws.Range("J85").Value = 35297.6
$ws.Range("L85").Value = 35297.6
$ws.Range("N85").Value = -37949.6

# Sheet BSM row 109
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

# Sheet BSM row 122
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 47420
$ws.Range("J122").Value = 47420
$ws.Range("L122").Value = 47420
$ws.Range("N122").Value = -57220

# Sheet BSM row 124
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents() | Out-Null

# Sheet BSM row 125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 54092
$ws.Range("J125").Value = 54092
$ws.Range("L125").Value = 54092
$ws.Range("N125").Value = -63932

# Sheet BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 31649.166
$ws.Range("J126").Value = 31649.166
$ws.Range("L126").Value = 31649.166
$ws.Range("N126").Value = -41529.166

# Sheet CRP row 109
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 25985
$ws.Range("J109").Value = 25985
$ws.Range("L109").Value = 25985
$ws.Range("N109").Value = -28065

# Sheet CRP row 127
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 54740
$ws.Range("J127").Value = 54740
$ws.Range("L127").Value = 54740
$ws.Range("N127").Value = -64660

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8643.137000000001
$ws.Range("I134").Value = 8643.137000000001
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 25929.411
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -23394.411
$ws.Range("N134").ClearContents() | Out-Null

# Sheet CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1264.5807
$ws.Range("I92").Value = 1157.2858
$ws.Range("J92").Value = 1352.9412
$ws.Range("K92").Value = 3471.8574
$ws.Range("L92").Value = 4058.8236
$ws.Range("M92").Value = -2223.8574
$ws.Range("N92").Value = -6554.8236

# Sheet GSM row 41
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 9800
$ws.Range("J41").Value = 9800
$ws.Range("L41").Value = 9800
$ws.Range("N41").Value = -10510

# Sheet GSM row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3342.6
$ws.Range("I43").Value = 886.8570999999999
$ws.Range("J43").Value = 9072.666999999999
$ws.Range("K43").Value = 886.8570999999999
$ws.Range("L43").Value = 9072.666999999999
$ws.Range("M43").Value = -735.8570999999999
$ws.Range("N43").Value = -9374.666999999999

# Sheet GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 13985.5
$ws.Range("J57").Value = 15473.917
$ws.Range("L57").Value = 15473.917
$ws.Range("N57").Value = -17113.917

# Sheet GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2200
$ws.Range("I80").Value = 2283.3333
$ws.Range("K80").Value = 2283.3333
$ws.Range("M80").Value = -1285.3333

# Sheet GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2200
$ws.Range("I83").Value = 2283.3333
$ws.Range("K83").Value = 11416.6665
$ws.Range("M83").Value = -6424.666499999999

# Sheet GSM row 93
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 11740.333
$ws.Range("J93").Value = 11740.333
$ws.Range("L93").Value = 11740.333
$ws.Range("N93").Value = -15484.333

# Sheet GSM row 111
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 22900
$ws.Range("J111").Value = 22900
$ws.Range("L111").Value = 22900
$ws.Range("N111").Value = -29034

# Sheet GSM row 124
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 33480
$ws.Range("J124").Value = 33480
$ws.Range("L124").Value = 33480
$ws.Range("N124").Value = -43300

# Sheet GSM row 128
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 48351.43
$ws.Range("J128").Value = 48351.43
$ws.Range("L128").Value = 48351.43
$ws.Range("N128").Value = -58311.43

# Sheet GSM row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 58423.81
$ws.Range("J135").Value = 58423.81
$ws.Range("L135").Value = 58423.81
$ws.Range("N135").Value = -68563.81

# Sheet LTW row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 23429.4
$ws.Range("J108").Value = 23429.4
$ws.Range("L108").Value = 23429.4
$ws.Range("N108").Value = -31109.4

# Sheet LTW row 109
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 24330
$ws.Range("J109").Value = 24330
$ws.Range("L109").Value = 24330
$ws.Range("N109").Value = -27104

# Sheet LTW row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 40398.168
$ws.Range("J123").Value = 40398.168
$ws.Range("L123").Value = 40398.168
$ws.Range("N123").Value = -50198.168

# Sheet LTW row 125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 54980
$ws.Range("J125").Value = 54980
$ws.Range("L125").Value = 54980
$ws.Range("N125").Value = -64820

# Sheet LTW row 130
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 42463.625
$ws.Range("J130").Value = 42463.625
$ws.Range("L130").Value = 42463.625
$ws.Range("N130").Value = -52503.625

# Sheet LTW row 134
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 100409.664
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 100409.664
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 100409.664
$ws.Range("M134").ClearContents() | Out-Null
$ws.Range("N134").Value = -110549.664

# Sheet WVR row 108
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 29960
$ws.Range("J108").Value = 29960
$ws.Range("L108").Value = 29960
$ws.Range("N108").Value = -37640

# Sheet WVR row 109
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 19577
$ws.Range("J109").Value = 19577
$ws.Range("L109").Value = 19577
$ws.Range("N109").Value = -22351

# Sheet WVR row 111
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 34944
$ws.Range("J111").Value = 34944
$ws.Range("L111").Value = 34944
$ws.Range("N111").Value = -43124

# Sheet WVR row 125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 31131.666
$ws.Range("J125").Value = 31131.666
$ws.Range("L125").Value = 31131.666
$ws.Range("N125").Value = -40971.666
